$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.614.74'
$ws.Range("E2").Value = '  -2.30%  '

$ws.Range("D3").Value = '1.843.33'
$ws.Range("E3").Value = '  -1.27%  '

$ws.Range("E4").Value = '  -0.16%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '314.37'
$ws.Range("E5").Value = '  -1.49%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9996'
$ws.Range("E6").Value = '  -0.12%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4243'
$ws.Range("E7").Value = '  -3.06%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3652'
$ws.Range("E8").Value = '  -1.56%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '45.77'
$ws.Range("E9").Value = '  +1.68%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.8988'
$ws.Range("E11").Value = '  -4.57%  '

$ws.Range("E12").Value = '  -3.53%  '

$ws.Range("D13").Value = '1.786.02'
$ws.Range("E13").Value = '  -4.42%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.387'
$ws.Range("E14").Value = '  -1.39%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.560'
$ws.Range("E15").Value = '  -2.62%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.06856'
$ws.Range("E16").Value = '  -0.15%  '

$ws.Range("E17").Value = '  -0.06%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '78.25'
$ws.Range("E18").Value = '  -5.11%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000008853'
$ws.Range("E19").Value = '  -2.93%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.9988'
$ws.Range("E20").Value = '  -0.22%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '15.57'
$ws.Range("E21").Value = '  -2.75%  '

$ws.Range("D22").Value = '27.611.41'
$ws.Range("E22").Value = '  -2.23%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.969'
$ws.Range("E23").Value = '  -3.30%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '10.57'
$ws.Range("E24").Value = '  -2.51%  '

$ws.Range("B25").Value = 'Toncoin'
$ws.Range("C25").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.040'
$ws.Range("E25").Value = '  +0.85%  '

$ws.Range("B26").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C26").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D26").Value = '1.987.49'
$ws.Range("E26").Value = '  -4.66%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '154.32'
$ws.Range("E27").Value = '  -0.37%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.28'
$ws.Range("E28").Value = '  -0.83%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.248'
$ws.Range("E29").Value = '  -1.70%  '

$ws.Range("B30").Value = 'LidoDAOToken'
$ws.Range("C30").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.828'
$ws.Range("E30").Value = '  +5.43%  '

$ws.Range("B31").Value = 'BitcoinCash'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '110.96'
$ws.Range("E31").Value = '  -2.82%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.08873'
$ws.Range("E32").Value = '  -1.88%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.7762'
$ws.Range("E33").Value = '  -3.40%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.554'
$ws.Range("E34").Value = '  -6.23%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.939'
$ws.Range("E35").Value = '  -0.71%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.092'
$ws.Range("E36").Value = '  -6.98%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.9988'
$ws.Range("E37").Value = '  -0.23%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.05428'
$ws.Range("E38").Value = '  -0.54%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.094'
$ws.Range("E39").Value = '  -2.34%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01926'
$ws.Range("E40").Value = '  -1.58%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.820'
$ws.Range("E41").Value = '  -4.93%  '

$ws.Range("E42").Value = '  -3.59%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.797'
$ws.Range("E43").Value = '  -4.86%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.1641'
$ws.Range("E44").Value = '  -2.13%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.239'
$ws.Range("E45").Value = '  -5.57%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.06637'
$ws.Range("E46").Value = '  -2.09%  '

$ws.Range("B47").Value = 'Decentraland'
$ws.Range("C47").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.4720'
$ws.Range("E47").Value = '  -3.14%  '

$ws.Range("B48").Value = 'Quant'
$ws.Range("C48").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '105.72'
$ws.Range("E48").Value = '  -2.11%  '

$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '10.33'
$ws.Range("E49").Value = '  -1.95%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.9988'
$ws.Range("E50").Value = '  -0.16%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.640'
$ws.Range("E51").Value = '  -2.52%  '
